# Apply the edits described by the diff:
#  - sheet "test_sheet1": add H3 = 0, H4 = "  c  ", C5 = "  a"
#  - update the active cell selection to H6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_sheet1")

# New cell values (order matters for shared-string table insertion order)
$ws.Range("H3").Value = 0
$ws.Range("C5").Value = "  a"
$ws.Range("H4").Value = "  c  "

# Update the selection shown in the sheet view
$ws.Activate()
$ws.Range("H6").Select()

$wb.Save()
